# Auto-generated edit script: restructure PlayerPerformance workbook
# NOTE: this runtime's PowerShell/COM bridge does not reliably propagate
# writes made through custom-function parameters bound to COM objects,
# so this script intentionally avoids helper functions and stays flat,
# operating directly on module-scope $wb / $sheetVar variables.
$wb = $excel.ActiveWorkbook

# --- Add 'Player Info' sheet at the front ---
$firstSheet = $wb.Worksheets.Item(1)
$playerInfoSheet = $wb.Worksheets.Add($firstSheet)
$playerInfoSheet.Name = "Player Info"
$nCols = 4
$headerData = New-Object 'object[,]' 1,$nCols
$headerData[0,0] = 'ID'
$headerData[0,1] = 'NAME'
$headerData[0,2] = 'BATTING_HAND'
$headerData[0,3] = 'BOWL_STYLE'
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(1,1), $playerInfoSheet.Cells.Item(1,$nCols)).Value = $headerData
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(1,1), $playerInfoSheet.Cells.Item(1,$nCols)).Font.Bold = $true
$nRows = 2
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(2,1), $playerInfoSheet.Cells.Item($nRows,1)).NumberFormat = "@"
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(2,2), $playerInfoSheet.Cells.Item($nRows,2)).NumberFormat = "@"
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(2,3), $playerInfoSheet.Cells.Item($nRows,3)).NumberFormat = "@"
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(2,4), $playerInfoSheet.Cells.Item($nRows,4)).NumberFormat = "@"
$arr = New-Object 'object[,]' 1,$nCols
$arr[0,0] = '3836'
$arr[0,1] = 'Andre Dwayne Russell'
$arr[0,2] = 'Right Handed'
$arr[0,3] = 'Right Arm Fast'
$playerInfoSheet.Range($playerInfoSheet.Cells.Item(2,1), $playerInfoSheet.Cells.Item($nRows,$nCols)).Value = $arr

# --- Rewrite 'ODI Batting' sheet data (MATCH_CARD_LINK -> MATCH_CODE) ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$nCols = 10
$headerData = New-Object 'object[,]' 1,$nCols
$headerData[0,0] = 'MATCH_NUMBER'
$headerData[0,1] = 'INNING_NUMBER'
$headerData[0,2] = 'MATCH_DATE'
$headerData[0,3] = 'MATCH_CODE'
$headerData[0,4] = 'MATCH_INNING'
$headerData[0,5] = 'OPPONENT'
$headerData[0,6] = 'VENUE'
$headerData[0,7] = 'DISMISSAL'
$headerData[0,8] = 'RUNS_SCORED'
$headerData[0,9] = 'BALLS_FACED'
$battingSheet.Range($battingSheet.Cells.Item(1,1), $battingSheet.Cells.Item(1,$nCols)).Value = $headerData
$battingSheet.Range($battingSheet.Cells.Item(1,1), $battingSheet.Cells.Item(1,$nCols)).Font.Bold = $true
$nRows = 57
$battingSheet.Range($battingSheet.Cells.Item(2,1), $battingSheet.Cells.Item($nRows,1)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,2), $battingSheet.Cells.Item($nRows,2)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,3), $battingSheet.Cells.Item($nRows,3)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,4), $battingSheet.Cells.Item($nRows,4)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,5), $battingSheet.Cells.Item($nRows,5)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,6), $battingSheet.Cells.Item($nRows,6)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,7), $battingSheet.Cells.Item($nRows,7)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,8), $battingSheet.Cells.Item($nRows,8)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,9), $battingSheet.Cells.Item($nRows,9)).NumberFormat = "@"
$battingSheet.Range($battingSheet.Cells.Item(2,10), $battingSheet.Cells.Item($nRows,10)).NumberFormat = "@"
$arr = New-Object 'object[,]' 56,$nCols
$arr[0,0] = '1'
$arr[0,1] = '1'
$arr[0,2] = '11/03/2011'
$arr[0,3] = '3261'
$arr[0,4] = '1st'
$arr[0,5] = 'Ireland'
$arr[0,6] = 'Punjab Cricket Association IS Bindra Stadium'
$arr[0,7] = 'b J F Mooney'
$arr[0,8] = '3'
$arr[0,9] = '7'
$arr[1,0] = '2'
$arr[1,1] = '2'
$arr[1,2] = '17/03/2011'
$arr[1,3] = '3270'
$arr[1,4] = '2nd'
$arr[1,5] = 'England'
$arr[1,6] = 'MA Chidambaram Stadium'
$arr[1,7] = 'lbw b J C Tredwell'
$arr[1,8] = '49'
$arr[1,9] = '46'
$arr[2,0] = '3'
$arr[2,1] = '3'
$arr[2,2] = '20/03/2011'
$arr[2,3] = '3276'
$arr[2,4] = '2nd'
$arr[2,5] = 'India'
$arr[2,6] = 'MA Chidambaram Stadium'
$arr[2,7] = 'c Y K Pathan b Yuvraj Singh'
$arr[2,8] = '0'
$arr[2,9] = '5'
$arr[3,0] = '4'
$arr[3,1] = $null
$arr[3,2] = '23/04/2011'
$arr[3,3] = '3287'
$arr[3,4] = '1st'
$arr[3,5] = 'Pakistan'
$arr[3,6] = 'Daren Sammy National Cricket Stadium'
$arr[3,7] = 'did not bat'
$arr[3,8] = '-'
$arr[3,9] = '-'
$arr[4,0] = '5'
$arr[4,1] = '4'
$arr[4,2] = '28/04/2011'
$arr[4,3] = '3289'
$arr[4,4] = '1st'
$arr[4,5] = 'Pakistan'
$arr[4,6] = 'Kensington Oval'
$arr[4,7] = 'st Mohammad Salman b Saeed Ajmal'
$arr[4,8] = '3'
$arr[4,9] = '9'
$arr[5,0] = '6'
$arr[5,1] = '5'
$arr[5,2] = '11/06/2011'
$arr[5,3] = '3296'
$arr[5,4] = '1st'
$arr[5,5] = 'India'
$arr[5,6] = 'Sir Vivian Richards Stadium'
$arr[5,7] = 'not out'
$arr[5,8] = '92*'
$arr[5,9] = '64'
$arr[6,0] = '7'
$arr[6,1] = '6'
$arr[6,2] = '13/06/2011'
$arr[6,3] = '3297'
$arr[6,4] = '1st'
$arr[6,5] = 'India'
$arr[6,6] = 'Sir Vivian Richards Stadium'
$arr[6,7] = 'c V Kohli b P Kumar'
$arr[6,8] = '25'
$arr[6,9] = '14'
$arr[7,0] = '8'
$arr[7,1] = $null
$arr[7,2] = '16/06/2011'
$arr[7,3] = '3298'
$arr[7,4] = '2nd'
$arr[7,5] = 'India'
$arr[7,6] = 'Sabina Park'
$arr[7,7] = 'did not bat'
$arr[7,8] = '-'
$arr[7,9] = '-'
$arr[8,0] = '9'
$arr[8,1] = '7'
$arr[8,2] = '13/10/2011'
$arr[8,3] = '3334'
$arr[8,4] = '1st'
$arr[8,5] = 'Bangladesh'
$arr[8,6] = 'Shere Bangla National Stadium'
$arr[8,7] = 'not out'
$arr[8,8] = '7*'
$arr[8,9] = '3'
$arr[9,0] = '10'
$arr[9,1] = $null
$arr[9,2] = '15/10/2011'
$arr[9,3] = '3336'
$arr[9,4] = '2nd'
$arr[9,5] = 'Bangladesh'
$arr[9,6] = 'Shere Bangla National Stadium'
$arr[9,7] = 'did not bat'
$arr[9,8] = '-'
$arr[9,9] = '-'
$arr[10,0] = '11'
$arr[10,1] = '8'
$arr[10,2] = '18/10/2011'
$arr[10,3] = '3338'
$arr[10,4] = '1st'
$arr[10,5] = 'Bangladesh'
$arr[10,6] = 'Zahur Ahmed Chowdhury Stadium'
$arr[10,7] = 'c & b Suhrawadi Shuvo'
$arr[10,8] = '2'
$arr[10,9] = '5'
$arr[11,0] = '12'
$arr[11,1] = '9'
$arr[11,2] = '29/11/2011'
$arr[11,3] = '3353'
$arr[11,4] = '1st'
$arr[11,5] = 'India'
$arr[11,6] = 'Barabati Stadium'
$arr[11,7] = 'b V R Aaron'
$arr[11,8] = '22'
$arr[11,9] = '20'
$arr[12,0] = '13'
$arr[12,1] = '10'
$arr[12,2] = '02/12/2011'
$arr[12,3] = '3355'
$arr[12,4] = '1st'
$arr[12,5] = 'India'
$arr[12,6] = 'Dr YS Rajasekhara Reddy Cricket Stadium'
$arr[12,7] = 'b R A Jadeja'
$arr[12,8] = '11'
$arr[12,9] = '13'
$arr[13,0] = '14'
$arr[13,1] = '11'
$arr[13,2] = '05/12/2011'
$arr[13,3] = '3357'
$arr[13,4] = '1st'
$arr[13,5] = 'India'
$arr[13,6] = 'Narendra Modi Stadium'
$arr[13,7] = 'not out'
$arr[13,8] = '40*'
$arr[13,9] = '18'
$arr[14,0] = '15'
$arr[14,1] = '12'
$arr[14,2] = '08/12/2011'
$arr[14,3] = '3359'
$arr[14,4] = '2nd'
$arr[14,5] = 'India'
$arr[14,6] = 'Holkar Cricket Stadium'
$arr[14,7] = 'st P A Patel b S K Raina'
$arr[14,8] = '29'
$arr[14,9] = '24'
$arr[15,0] = '16'
$arr[15,1] = '13'
$arr[15,2] = '11/12/2011'
$arr[15,3] = '3360'
$arr[15,4] = '2nd'
$arr[15,5] = 'India'
$arr[15,6] = 'MA Chidambaram Stadium'
$arr[15,7] = 'run out'
$arr[15,8] = '53'
$arr[15,9] = '42'
$arr[16,0] = '17'
$arr[16,1] = '14'
$arr[16,2] = '16/03/2012'
$arr[16,3] = '3398'
$arr[16,4] = '2nd'
$arr[16,5] = 'Australia'
$arr[16,6] = 'Arnos Vale Ground'
$arr[16,7] = 'st M S Wade b X J Doherty'
$arr[16,8] = '1'
$arr[16,9] = '4'
$arr[17,0] = '18'
$arr[17,1] = $null
$arr[17,2] = '18/03/2012'
$arr[17,3] = '3400'
$arr[17,4] = '2nd'
$arr[17,5] = 'Australia'
$arr[17,6] = 'Arnos Vale Ground'
$arr[17,7] = 'did not bat'
$arr[17,8] = '-'
$arr[17,9] = '-'
$arr[18,0] = '19'
$arr[18,1] = '15'
$arr[18,2] = '20/03/2012'
$arr[18,3] = '3402'
$arr[18,4] = '2nd'
$arr[18,5] = 'Australia'
$arr[18,6] = 'Arnos Vale Ground'
$arr[18,7] = 'c †M S Wade b C J McKay'
$arr[18,8] = '37'
$arr[18,9] = '42'
$arr[19,0] = '20'
$arr[19,1] = '16'
$arr[19,2] = '23/03/2012'
$arr[19,3] = '3404'
$arr[19,4] = '1st'
$arr[19,5] = 'Australia'
$arr[19,6] = 'Daren Sammy National Cricket Stadium'
$arr[19,7] = 'c G J Bailey b S R Watson'
$arr[19,8] = '34'
$arr[19,9] = '32'
$arr[20,0] = '21'
$arr[20,1] = '17'
$arr[20,2] = '25/03/2012'
$arr[20,3] = '3405'
$arr[20,4] = '2nd'
$arr[20,5] = 'Australia'
$arr[20,6] = 'Daren Sammy National Cricket Stadium'
$arr[20,7] = 'lbw b X J Doherty'
$arr[20,8] = '41'
$arr[20,9] = '33'
$arr[21,0] = '22'
$arr[21,1] = '18'
$arr[21,2] = '16/06/2012'
$arr[21,3] = '3412'
$arr[21,4] = '2nd'
$arr[21,5] = 'England'
$arr[21,6] = 'The Rose Bowl'
$arr[21,7] = 'c E J G Morgan b T T Bresnan'
$arr[21,8] = '7'
$arr[21,9] = '13'
$arr[22,0] = '23'
$arr[22,1] = $null
$arr[22,2] = '05/07/2012'
$arr[22,3] = '3422'
$arr[22,4] = '2nd'
$arr[22,5] = 'New Zealand'
$arr[22,6] = 'Sabina Park'
$arr[22,7] = 'did not bat'
$arr[22,8] = '-'
$arr[22,9] = '-'
$arr[23,0] = '24'
$arr[23,1] = $null
$arr[23,2] = '07/07/2012'
$arr[23,3] = '3424'
$arr[23,4] = '1st'
$arr[23,5] = 'New Zealand'
$arr[23,6] = 'Sabina Park'
$arr[23,7] = 'did not bat'
$arr[23,8] = '-'
$arr[23,9] = '-'
$arr[24,0] = '25'
$arr[24,1] = '19'
$arr[24,2] = '11/07/2012'
$arr[24,3] = '3427'
$arr[24,4] = '2nd'
$arr[24,5] = 'New Zealand'
$arr[24,6] = 'Warner Park'
$arr[24,7] = 'not out'
$arr[24,8] = '42*'
$arr[24,9] = '24'
$arr[25,0] = '26'
$arr[25,1] = '20'
$arr[25,2] = '14/07/2012'
$arr[25,3] = '3429'
$arr[25,4] = '1st'
$arr[25,5] = 'New Zealand'
$arr[25,6] = 'Warner Park'
$arr[25,7] = 'c M J Guptill b J D P Oram'
$arr[25,8] = '29'
$arr[25,9] = '16'
$arr[26,0] = '27'
$arr[26,1] = '21'
$arr[26,2] = '16/07/2012'
$arr[26,3] = '3430'
$arr[26,4] = '1st'
$arr[26,5] = 'New Zealand'
$arr[26,6] = 'Warner Park'
$arr[26,7] = 'not out'
$arr[26,8] = '59*'
$arr[26,9] = '40'
$arr[27,0] = '28'
$arr[27,1] = '22'
$arr[27,2] = '30/11/2012'
$arr[27,3] = '3450'
$arr[27,4] = '1st'
$arr[27,5] = 'Bangladesh'
$arr[27,6] = 'Sheikh Abu Naser Stadium'
$arr[27,7] = 'b Abdur Razzak'
$arr[27,8] = '0'
$arr[27,9] = '3'
$arr[28,0] = '29'
$arr[28,1] = '23'
$arr[28,2] = '02/12/2012'
$arr[28,3] = '3451'
$arr[28,4] = '2nd'
$arr[28,5] = 'Bangladesh'
$arr[28,6] = 'Sheikh Abu Naser Stadium'
$arr[28,7] = 'run out'
$arr[28,8] = '9'
$arr[28,9] = '15'
$arr[29,0] = '30'
$arr[29,1] = '24'
$arr[29,2] = '08/12/2012'
$arr[29,3] = '3454'
$arr[29,4] = '1st'
$arr[29,5] = 'Bangladesh'
$arr[29,6] = 'Shere Bangla National Stadium'
$arr[29,7] = 'lbw b Mahmudullah'
$arr[29,8] = '0'
$arr[29,9] = '4'
$arr[30,0] = '31'
$arr[30,1] = '25'
$arr[30,2] = '06/02/2013'
$arr[30,3] = '3473'
$arr[30,4] = '2nd'
$arr[30,5] = 'Australia'
$arr[30,6] = 'Manuka Oval'
$arr[30,7] = 'c †M S Wade b C J McKay'
$arr[30,8] = '43'
$arr[30,9] = '31'
$arr[31,0] = '32'
$arr[31,1] = '26'
$arr[31,2] = '08/02/2013'
$arr[31,3] = '3474'
$arr[31,4] = '1st'
$arr[31,5] = 'Australia'
$arr[31,6] = 'Sydney Cricket Ground'
$arr[31,7] = 'c G J Maxwell b J P Faulkner'
$arr[31,8] = '18'
$arr[31,9] = '17'
$arr[32,0] = '33'
$arr[32,1] = '27'
$arr[32,2] = '22/02/2013'
$arr[32,3] = '3478'
$arr[32,4] = '1st'
$arr[32,5] = 'Zimbabwe'
$arr[32,6] = 'National Cricket Stadium (Grenada)'
$arr[32,7] = 'c R W Chakabva b N M''shangwe'
$arr[32,8] = '4'
$arr[32,9] = '9'
$arr[33,0] = '34'
$arr[33,1] = $null
$arr[33,2] = '24/02/2013'
$arr[33,3] = '3480'
$arr[33,4] = '2nd'
$arr[33,5] = 'Zimbabwe'
$arr[33,6] = 'National Cricket Stadium (Grenada)'
$arr[33,7] = 'did not bat'
$arr[33,8] = '-'
$arr[33,9] = '-'
$arr[34,0] = '35'
$arr[34,1] = '28'
$arr[34,2] = '07/01/2014'
$arr[34,3] = '3598'
$arr[34,4] = '1st'
$arr[34,5] = 'New Zealand'
$arr[34,6] = 'Seddon Park'
$arr[34,7] = 'not out'
$arr[34,8] = '6*'
$arr[34,9] = '3'
$arr[35,0] = '36'
$arr[35,1] = '29'
$arr[35,2] = '08/10/2014'
$arr[35,3] = '3678'
$arr[35,4] = '1st'
$arr[35,5] = 'India'
$arr[35,6] = 'Nehru Stadium (Kochi)'
$arr[35,7] = 'c V Kohli b Mohammed Shami'
$arr[35,8] = '1'
$arr[35,9] = '2'
$arr[36,0] = '37'
$arr[36,1] = '30'
$arr[36,2] = '11/10/2014'
$arr[36,3] = '3680'
$arr[36,4] = '2nd'
$arr[36,5] = 'India'
$arr[36,6] = 'Arun Jaitley Stadium'
$arr[36,7] = 'st M S Dhoni b R A Jadeja'
$arr[36,8] = '4'
$arr[36,9] = '6'
$arr[37,0] = '38'
$arr[37,1] = '31'
$arr[37,2] = '17/10/2014'
$arr[37,3] = '3683'
$arr[37,4] = '2nd'
$arr[37,5] = 'India'
$arr[37,6] = 'Himachal Pradesh Cricket Association Stadium'
$arr[37,7] = 'b U T Yadav'
$arr[37,8] = '46'
$arr[37,9] = '23'
$arr[38,0] = '39'
$arr[38,1] = '32'
$arr[38,2] = '16/01/2015'
$arr[38,3] = '3728'
$arr[38,4] = '2nd'
$arr[38,5] = 'South Africa'
$arr[38,6] = 'Kingsmead'
$arr[38,7] = 'c †A B de Villiers b V D Philander'
$arr[38,8] = '19'
$arr[38,9] = '8'
$arr[39,0] = '40'
$arr[39,1] = '33'
$arr[39,2] = '18/01/2015'
$arr[39,3] = '3732'
$arr[39,4] = '2nd'
$arr[39,5] = 'South Africa'
$arr[39,6] = 'Wanderers Stadium'
$arr[39,7] = 'c F du Plessis b Imran Tahir'
$arr[39,8] = '0'
$arr[39,9] = '1'
$arr[40,0] = '41'
$arr[40,1] = '34'
$arr[40,2] = '21/01/2015'
$arr[40,3] = '3736'
$arr[40,4] = '1st'
$arr[40,5] = 'South Africa'
$arr[40,6] = 'Buffalo Park'
$arr[40,7] = 'b M Morkel'
$arr[40,8] = '16'
$arr[40,9] = '25'
$arr[41,0] = '42'
$arr[41,1] = '35'
$arr[41,2] = '25/01/2015'
$arr[41,3] = '3740'
$arr[41,4] = '2nd'
$arr[41,5] = 'South Africa'
$arr[41,6] = 'St George''s Park'
$arr[41,7] = 'not out'
$arr[41,8] = '64*'
$arr[41,9] = '40'
$arr[42,0] = '43'
$arr[42,1] = '36'
$arr[42,2] = '28/01/2015'
$arr[42,3] = '3742'
$arr[42,4] = '2nd'
$arr[42,5] = 'South Africa'
$arr[42,6] = 'SuperSport Park'
$arr[42,7] = 'c †Q de Kock b W D Parnell'
$arr[42,8] = '24'
$arr[42,9] = '12'
$arr[43,0] = '44'
$arr[43,1] = '37'
$arr[43,2] = '16/02/2015'
$arr[43,3] = '3752'
$arr[43,4] = '1st'
$arr[43,5] = 'Ireland'
$arr[43,6] = 'Saxton Oval'
$arr[43,7] = 'not out'
$arr[43,8] = '27*'
$arr[43,9] = '13'
$arr[44,0] = '45'
$arr[44,1] = '38'
$arr[44,2] = '21/02/2015'
$arr[44,3] = '3757'
$arr[44,4] = '1st'
$arr[44,5] = 'Pakistan'
$arr[44,6] = 'Hagley Oval'
$arr[44,7] = 'not out'
$arr[44,8] = '42*'
$arr[44,9] = '13'
$arr[45,0] = '46'
$arr[45,1] = $null
$arr[45,2] = '24/02/2015'
$arr[45,3] = '3762'
$arr[45,4] = '1st'
$arr[45,5] = 'Zimbabwe'
$arr[45,6] = 'Manuka Oval'
$arr[45,7] = 'did not bat'
$arr[45,8] = '-'
$arr[45,9] = '-'
$arr[46,0] = '47'
$arr[46,1] = '39'
$arr[46,2] = '27/02/2015'
$arr[46,3] = '3766'
$arr[46,4] = '2nd'
$arr[46,5] = 'South Africa'
$arr[46,6] = 'Sydney Cricket Ground'
$arr[46,7] = 'c K J Abbott b Imran Tahir'
$arr[46,8] = '0'
$arr[46,9] = '3'
$arr[47,0] = '48'
$arr[47,1] = '40'
$arr[47,2] = '06/03/2015'
$arr[47,3] = '3775'
$arr[47,4] = '1st'
$arr[47,5] = 'India'
$arr[47,6] = 'WACA Ground'
$arr[47,7] = 'c V Kohli b R A Jadeja'
$arr[47,8] = '8'
$arr[47,9] = '8'
$arr[48,0] = '49'
$arr[48,1] = '41'
$arr[48,2] = '15/03/2015'
$arr[48,3] = '3788'
$arr[48,4] = '2nd'
$arr[48,5] = 'United Arab Emirates'
$arr[48,6] = 'McLean Park'
$arr[48,7] = 'c & b Amjad Javed'
$arr[48,8] = '7'
$arr[48,9] = '8'
$arr[49,0] = '50'
$arr[49,1] = '42'
$arr[49,2] = '21/03/2015'
$arr[49,3] = '3793'
$arr[49,4] = '2nd'
$arr[49,5] = 'New Zealand'
$arr[49,6] = 'Sky Stadium'
$arr[49,7] = 'b T G Southee'
$arr[49,8] = '20'
$arr[49,9] = '11'
$arr[50,0] = '51'
$arr[50,1] = '43'
$arr[50,2] = '01/11/2015'
$arr[50,3] = '3852'
$arr[50,4] = '1st'
$arr[50,5] = 'Sri Lanka'
$arr[50,6] = 'R Premadasa Stadium'
$arr[50,7] = 'c M D Gunathilaka b A D Mathews'
$arr[50,8] = '41'
$arr[50,9] = '24'
$arr[51,0] = '52'
$arr[51,1] = '44'
$arr[51,2] = '22/07/2018'
$arr[51,3] = '4179'
$arr[51,4] = '2nd'
$arr[51,5] = 'Bangladesh'
$arr[51,6] = 'Providence Stadium'
$arr[51,7] = 'c Mahmudullah b Mashrafe Mortaza'
$arr[51,8] = '13'
$arr[51,9] = '12'
$arr[52,0] = '53'
$arr[52,1] = $null
$arr[52,2] = '31/05/2019'
$arr[52,3] = '4304'
$arr[52,4] = '2nd'
$arr[52,5] = 'Pakistan'
$arr[52,6] = 'Trent Bridge'
$arr[52,7] = 'did not bat'
$arr[52,8] = '-'
$arr[52,9] = '-'
$arr[53,0] = '54'
$arr[53,1] = '45'
$arr[53,2] = '06/06/2019'
$arr[53,3] = '4312'
$arr[53,4] = '2nd'
$arr[53,5] = 'Australia'
$arr[53,6] = 'Trent Bridge'
$arr[53,7] = 'c G J Maxwell b M A Starc'
$arr[53,8] = '15'
$arr[53,9] = '11'
$arr[54,0] = '55'
$arr[54,1] = '46'
$arr[54,2] = '14/06/2019'
$arr[54,3] = '4321'
$arr[54,4] = '1st'
$arr[54,5] = 'England'
$arr[54,6] = 'The Rose Bowl'
$arr[54,7] = 'c C R Woakes b M A Wood'
$arr[54,8] = '21'
$arr[54,9] = '16'
$arr[55,0] = '56'
$arr[55,1] = '47'
$arr[55,2] = '17/06/2019'
$arr[55,3] = '4325'
$arr[55,4] = '1st'
$arr[55,5] = 'Bangladesh'
$arr[55,6] = 'The Cooper Associates County Ground'
$arr[55,7] = 'c †Mushfiqur Rahim b Mustafizur Rahman'
$arr[55,8] = '0'
$arr[55,9] = '2'
$battingSheet.Range($battingSheet.Cells.Item(2,1), $battingSheet.Cells.Item($nRows,$nCols)).Value = $arr

# --- Rewrite 'ODI Bowling' sheet data (MATCH_CARD_LINK -> MATCH_CODE) ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$nCols = 7
$headerData = New-Object 'object[,]' 1,$nCols
$headerData[0,0] = 'MATCH_NUMBER'
$headerData[0,1] = 'MATCH_CODE'
$headerData[0,2] = 'MATCH_INNING'
$headerData[0,3] = 'OPPONENT'
$headerData[0,4] = 'VENUE'
$headerData[0,5] = 'OVERS'
$headerData[0,6] = 'WICKETS_RUNS'
$bowlingSheet.Range($bowlingSheet.Cells.Item(1,1), $bowlingSheet.Cells.Item(1,$nCols)).Value = $headerData
$bowlingSheet.Range($bowlingSheet.Cells.Item(1,1), $bowlingSheet.Cells.Item(1,$nCols)).Font.Bold = $true
$nRows = 56
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,1), $bowlingSheet.Cells.Item($nRows,1)).NumberFormat = "@"
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,2), $bowlingSheet.Cells.Item($nRows,2)).NumberFormat = "@"
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,3), $bowlingSheet.Cells.Item($nRows,3)).NumberFormat = "@"
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,4), $bowlingSheet.Cells.Item($nRows,4)).NumberFormat = "@"
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,5), $bowlingSheet.Cells.Item($nRows,5)).NumberFormat = "@"
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,6), $bowlingSheet.Cells.Item($nRows,6)).NumberFormat = "@"
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,7), $bowlingSheet.Cells.Item($nRows,7)).NumberFormat = "@"
$arr = New-Object 'object[,]' 55,$nCols
$arr[0,0] = '1'
$arr[0,1] = '3261'
$arr[0,2] = '2nd'
$arr[0,3] = 'Ireland'
$arr[0,4] = 'Punjab Cricket Association IS Bindra Stadium'
$arr[0,5] = '10.0'
$arr[0,6] = '1/37'
$arr[1,0] = '2'
$arr[1,1] = '3270'
$arr[1,2] = '1st'
$arr[1,3] = 'England'
$arr[1,4] = 'MA Chidambaram Stadium'
$arr[1,5] = '8.0'
$arr[1,6] = '4/49'
$arr[2,0] = '3'
$arr[2,1] = '3276'
$arr[2,2] = '1st'
$arr[2,3] = 'India'
$arr[2,4] = 'MA Chidambaram Stadium'
$arr[2,5] = '9.1'
$arr[2,6] = '2/46'
$arr[3,0] = '4'
$arr[3,1] = '3287'
$arr[3,2] = '2nd'
$arr[3,3] = 'Pakistan'
$arr[3,4] = 'Daren Sammy National Cricket Stadium'
$arr[3,5] = '5.0'
$arr[3,6] = '0/36'
$arr[4,0] = '5'
$arr[4,1] = '3289'
$arr[4,2] = '2nd'
$arr[4,3] = 'Pakistan'
$arr[4,4] = 'Kensington Oval'
$arr[4,5] = '4.0'
$arr[4,6] = '0/23'
$arr[5,0] = '6'
$arr[5,1] = '3296'
$arr[5,2] = '2nd'
$arr[5,3] = 'India'
$arr[5,4] = 'Sir Vivian Richards Stadium'
$arr[5,5] = '9.0'
$arr[5,6] = '1/59'
$arr[6,0] = '7'
$arr[6,1] = '3297'
$arr[6,2] = '2nd'
$arr[6,3] = 'India'
$arr[6,4] = 'Sir Vivian Richards Stadium'
$arr[6,5] = '7.0'
$arr[6,6] = '3/16'
$arr[7,0] = '8'
$arr[7,1] = '3298'
$arr[7,2] = '1st'
$arr[7,3] = 'India'
$arr[7,4] = 'Sabina Park'
$arr[7,5] = '8.3'
$arr[7,6] = '4/35'
$arr[8,0] = '9'
$arr[8,1] = '3334'
$arr[8,2] = '2nd'
$arr[8,3] = 'Bangladesh'
$arr[8,4] = 'Shere Bangla National Stadium'
$arr[8,5] = '9.0'
$arr[8,6] = '2/44'
$arr[9,0] = '10'
$arr[9,1] = '3336'
$arr[9,2] = '1st'
$arr[9,3] = 'Bangladesh'
$arr[9,4] = 'Shere Bangla National Stadium'
$arr[9,5] = '9.0'
$arr[9,6] = '1/50'
$arr[10,0] = '12'
$arr[10,1] = '3353'
$arr[10,2] = '2nd'
$arr[10,3] = 'India'
$arr[10,4] = 'Barabati Stadium'
$arr[10,5] = '9.0'
$arr[10,6] = '2/29'
$arr[11,0] = '13'
$arr[11,1] = '3355'
$arr[11,2] = '2nd'
$arr[11,3] = 'India'
$arr[11,4] = 'Dr YS Rajasekhara Reddy Cricket Stadium'
$arr[11,5] = '8.1'
$arr[11,6] = '0/60'
$arr[12,0] = '14'
$arr[12,1] = '3357'
$arr[12,2] = '2nd'
$arr[12,3] = 'India'
$arr[12,4] = 'Narendra Modi Stadium'
$arr[12,5] = '4.0'
$arr[12,6] = '0/25'
$arr[13,0] = '15'
$arr[13,1] = '3359'
$arr[13,2] = '1st'
$arr[13,3] = 'India'
$arr[13,4] = 'Holkar Cricket Stadium'
$arr[13,5] = '7.0'
$arr[13,6] = '1/63'
$arr[14,0] = '16'
$arr[14,1] = '3360'
$arr[14,2] = '1st'
$arr[14,3] = 'India'
$arr[14,4] = 'MA Chidambaram Stadium'
$arr[14,5] = '6.0'
$arr[14,6] = '0/31'
$arr[15,0] = '17'
$arr[15,1] = '3398'
$arr[15,2] = '1st'
$arr[15,3] = 'Australia'
$arr[15,4] = 'Arnos Vale Ground'
$arr[15,5] = '4.0'
$arr[15,6] = '0/21'
$arr[16,0] = '18'
$arr[16,1] = '3400'
$arr[16,2] = '1st'
$arr[16,3] = 'Australia'
$arr[16,4] = 'Arnos Vale Ground'
$arr[16,5] = '6.0'
$arr[16,6] = '0/14'
$arr[17,0] = '19'
$arr[17,1] = '3402'
$arr[17,2] = '1st'
$arr[17,3] = 'Australia'
$arr[17,4] = 'Arnos Vale Ground'
$arr[17,5] = '5.0'
$arr[17,6] = '0/28'
$arr[18,0] = '20'
$arr[18,1] = '3404'
$arr[18,2] = '2nd'
$arr[18,3] = 'Australia'
$arr[18,4] = 'Daren Sammy National Cricket Stadium'
$arr[18,5] = '7.0'
$arr[18,6] = '2/34'
$arr[19,0] = '21'
$arr[19,1] = '3405'
$arr[19,2] = '1st'
$arr[19,3] = 'Australia'
$arr[19,4] = 'Daren Sammy National Cricket Stadium'
$arr[19,5] = '9.0'
$arr[19,6] = '4/61'
$arr[20,0] = '22'
$arr[20,1] = '3412'
$arr[20,2] = '1st'
$arr[20,3] = 'England'
$arr[20,4] = 'The Rose Bowl'
$arr[20,5] = '6.0'
$arr[20,6] = '0/43'
$arr[21,0] = '23'
$arr[21,1] = '3422'
$arr[21,2] = '1st'
$arr[21,3] = 'New Zealand'
$arr[21,4] = 'Sabina Park'
$arr[21,5] = '10.0'
$arr[21,6] = '4/45'
$arr[22,0] = '24'
$arr[22,1] = '3424'
$arr[22,2] = '2nd'
$arr[22,3] = 'New Zealand'
$arr[22,4] = 'Sabina Park'
$arr[22,5] = '9.0'
$arr[22,6] = '1/51'
$arr[23,0] = '25'
$arr[23,1] = '3427'
$arr[23,2] = '1st'
$arr[23,3] = 'New Zealand'
$arr[23,4] = 'Warner Park'
$arr[23,5] = '9.0'
$arr[23,6] = '4/57'
$arr[24,0] = '26'
$arr[24,1] = '3429'
$arr[24,2] = '2nd'
$arr[24,3] = 'New Zealand'
$arr[24,4] = 'Warner Park'
$arr[24,5] = '6.0'
$arr[24,6] = '1/49'
$arr[25,0] = '27'
$arr[25,1] = '3430'
$arr[25,2] = '2nd'
$arr[25,3] = 'New Zealand'
$arr[25,4] = 'Warner Park'
$arr[25,5] = '4.0'
$arr[25,6] = '0/25'
$arr[26,0] = '28'
$arr[26,1] = '3450'
$arr[26,2] = '2nd'
$arr[26,3] = 'Bangladesh'
$arr[26,4] = 'Sheikh Abu Naser Stadium'
$arr[26,5] = '7.0'
$arr[26,6] = '1/28'
$arr[27,0] = '29'
$arr[27,1] = '3451'
$arr[27,2] = '1st'
$arr[27,3] = 'Bangladesh'
$arr[27,4] = 'Sheikh Abu Naser Stadium'
$arr[27,5] = '9.0'
$arr[27,6] = '1/58'
$arr[28,0] = '30'
$arr[28,1] = '3454'
$arr[28,2] = '2nd'
$arr[28,3] = 'Bangladesh'
$arr[28,4] = 'Shere Bangla National Stadium'
$arr[28,5] = '10.0'
$arr[28,6] = '0/51'
$arr[29,0] = '31'
$arr[29,1] = '3473'
$arr[29,2] = '1st'
$arr[29,3] = 'Australia'
$arr[29,4] = 'Manuka Oval'
$arr[29,5] = '5.0'
$arr[29,6] = '0/41'
$arr[30,0] = '32'
$arr[30,1] = '3474'
$arr[30,2] = '2nd'
$arr[30,3] = 'Australia'
$arr[30,4] = 'Sydney Cricket Ground'
$arr[30,5] = '6.0'
$arr[30,6] = '1/40'
$arr[31,0] = '33'
$arr[31,1] = '3478'
$arr[31,2] = '2nd'
$arr[31,3] = 'Zimbabwe'
$arr[31,4] = 'National Cricket Stadium (Grenada)'
$arr[31,5] = '7.0'
$arr[31,6] = '2/24'
$arr[32,0] = '34'
$arr[32,1] = '3480'
$arr[32,2] = '1st'
$arr[32,3] = 'Zimbabwe'
$arr[32,4] = 'National Cricket Stadium (Grenada)'
$arr[32,5] = '6.0'
$arr[32,6] = '0/29'
$arr[33,0] = '35'
$arr[33,1] = '3598'
$arr[33,2] = '2nd'
$arr[33,3] = 'New Zealand'
$arr[33,4] = 'Seddon Park'
$arr[33,5] = '4.0'
$arr[33,6] = '2/31'
$arr[34,0] = '36'
$arr[34,1] = '3678'
$arr[34,2] = '2nd'
$arr[34,3] = 'India'
$arr[34,4] = 'Nehru Stadium (Kochi)'
$arr[34,5] = '4.0'
$arr[34,6] = '1/21'
$arr[35,0] = '37'
$arr[35,1] = '3680'
$arr[35,2] = '1st'
$arr[35,3] = 'India'
$arr[35,4] = 'Arun Jaitley Stadium'
$arr[35,5] = '3.0'
$arr[35,6] = '0/14'
$arr[36,0] = '38'
$arr[36,1] = '3683'
$arr[36,2] = '1st'
$arr[36,3] = 'India'
$arr[36,4] = 'Himachal Pradesh Cricket Association Stadium'
$arr[36,5] = '7.0'
$arr[36,6] = '1/48'
$arr[37,0] = '39'
$arr[37,1] = '3728'
$arr[37,2] = '1st'
$arr[37,3] = 'South Africa'
$arr[37,4] = 'Kingsmead'
$arr[37,5] = '10.0'
$arr[37,6] = '2/51'
$arr[38,0] = '40'
$arr[38,1] = '3732'
$arr[38,2] = '1st'
$arr[38,3] = 'South Africa'
$arr[38,4] = 'Wanderers Stadium'
$arr[38,5] = '10.0'
$arr[38,6] = '1/78'
$arr[39,0] = '41'
$arr[39,1] = '3736'
$arr[39,2] = '2nd'
$arr[39,3] = 'South Africa'
$arr[39,4] = 'Buffalo Park'
$arr[39,5] = '5.0'
$arr[39,6] = '0/15'
$arr[40,0] = '42'
$arr[40,1] = '3740'
$arr[40,2] = '1st'
$arr[40,3] = 'South Africa'
$arr[40,4] = 'St George''s Park'
$arr[40,5] = '10.0'
$arr[40,6] = '1/60'
$arr[41,0] = '43'
$arr[41,1] = '3742'
$arr[41,2] = '1st'
$arr[41,3] = 'South Africa'
$arr[41,4] = 'SuperSport Park'
$arr[41,5] = '8.0'
$arr[41,6] = '3/85'
$arr[42,0] = '44'
$arr[42,1] = '3752'
$arr[42,2] = '2nd'
$arr[42,3] = 'Ireland'
$arr[42,4] = 'Saxton Oval'
$arr[42,5] = '6.0'
$arr[42,6] = '0/33'
$arr[43,0] = '45'
$arr[43,1] = '3757'
$arr[43,2] = '2nd'
$arr[43,3] = 'Pakistan'
$arr[43,4] = 'Hagley Oval'
$arr[43,5] = '8.0'
$arr[43,6] = '3/33'
$arr[44,0] = '46'
$arr[44,1] = '3762'
$arr[44,2] = '2nd'
$arr[44,3] = 'Zimbabwe'
$arr[44,4] = 'Manuka Oval'
$arr[44,5] = '5.0'
$arr[44,6] = '0/44'
$arr[45,0] = '47'
$arr[45,1] = '3766'
$arr[45,2] = '1st'
$arr[45,3] = 'South Africa'
$arr[45,4] = 'Sydney Cricket Ground'
$arr[45,5] = '9.0'
$arr[45,6] = '2/74'
$arr[46,0] = '48'
$arr[46,1] = '3775'
$arr[46,2] = '2nd'
$arr[46,3] = 'India'
$arr[46,4] = 'WACA Ground'
$arr[46,5] = '8.0'
$arr[46,6] = '2/43'
$arr[47,0] = '49'
$arr[47,1] = '3788'
$arr[47,2] = '1st'
$arr[47,3] = 'United Arab Emirates'
$arr[47,4] = 'McLean Park'
$arr[47,5] = '8.0'
$arr[47,6] = '2/20'
$arr[48,0] = '50'
$arr[48,1] = '3793'
$arr[48,2] = '1st'
$arr[48,3] = 'New Zealand'
$arr[48,4] = 'Sky Stadium'
$arr[48,5] = '10.0'
$arr[48,6] = '2/96'
$arr[49,0] = '51'
$arr[49,1] = '3852'
$arr[49,2] = '2nd'
$arr[49,3] = 'Sri Lanka'
$arr[49,4] = 'R Premadasa Stadium'
$arr[49,5] = '0.5'
$arr[49,6] = '0/18'
$arr[50,0] = '52'
$arr[50,1] = '4179'
$arr[50,2] = '1st'
$arr[50,3] = 'Bangladesh'
$arr[50,4] = 'Providence Stadium'
$arr[50,5] = '9.0'
$arr[50,6] = '1/62'
$arr[51,0] = '53'
$arr[51,1] = '4304'
$arr[51,2] = '1st'
$arr[51,3] = 'Pakistan'
$arr[51,4] = 'Trent Bridge'
$arr[51,5] = '3.0'
$arr[51,6] = '2/4'
$arr[52,0] = '54'
$arr[52,1] = '4312'
$arr[52,2] = '1st'
$arr[52,3] = 'Australia'
$arr[52,4] = 'Trent Bridge'
$arr[52,5] = '8.0'
$arr[52,6] = '2/41'
$arr[53,0] = '55'
$arr[53,1] = '4321'
$arr[53,2] = '2nd'
$arr[53,3] = 'England'
$arr[53,4] = 'The Rose Bowl'
$arr[53,5] = '2.0'
$arr[53,6] = '0/14'
$arr[54,0] = '56'
$arr[54,1] = '4325'
$arr[54,2] = '2nd'
$arr[54,3] = 'Bangladesh'
$arr[54,4] = 'The Cooper Associates County Ground'
$arr[54,5] = '6.0'
$arr[54,6] = '1/42'
$bowlingSheet.Range($bowlingSheet.Cells.Item(2,1), $bowlingSheet.Cells.Item($nRows,$nCols)).Value = $arr

# --- Add 'ODI Batting Extra' sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add($null, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"
$nCols = 6
$headerData = New-Object 'object[,]' 1,$nCols
$headerData[0,0] = 'MATCH_CODE'
$headerData[0,1] = 'BATTING_POSITION'
$headerData[0,2] = 'NUM_4'
$headerData[0,3] = 'NUM_6'
$headerData[0,4] = 'PERCENT_RUNS_OF_TOTAL'
$headerData[0,5] = 'MAN_OF_MATCH'
$extraSheet.Range($extraSheet.Cells.Item(1,1), $extraSheet.Cells.Item(1,$nCols)).Value = $headerData
$extraSheet.Range($extraSheet.Cells.Item(1,1), $extraSheet.Cells.Item(1,$nCols)).Font.Bold = $true
$nRows = 21
$extraSheet.Range($extraSheet.Cells.Item(2,1), $extraSheet.Cells.Item($nRows,1)).NumberFormat = "@"
$extraSheet.Range($extraSheet.Cells.Item(2,3), $extraSheet.Cells.Item($nRows,3)).NumberFormat = "@"
$extraSheet.Range($extraSheet.Cells.Item(2,4), $extraSheet.Cells.Item($nRows,4)).NumberFormat = "@"
$extraSheet.Range($extraSheet.Cells.Item(2,5), $extraSheet.Cells.Item($nRows,5)).NumberFormat = "@"
$extraSheet.Range($extraSheet.Cells.Item(2,6), $extraSheet.Cells.Item($nRows,6)).NumberFormat = "@"
$arr = New-Object 'object[,]' 20,$nCols
$arr[0,0] = '3680'
$arr[0,1] = 7
$arr[0,2] = '0'
$arr[0,3] = '0'
$arr[0,4] = '1.86%'
$arr[0,5] = 'NO'
$arr[1,0] = '3683'
$arr[1,1] = 8
$arr[1,2] = '6'
$arr[1,3] = '3'
$arr[1,4] = '16.97%'
$arr[1,5] = 'NO'
$arr[2,0] = '3728'
$arr[2,1] = $null
$arr[2,2] = $null
$arr[2,3] = $null
$arr[2,4] = $null
$arr[2,5] = 'NO'
$arr[3,0] = '3732'
$arr[3,1] = 7
$arr[3,2] = '0'
$arr[3,3] = '0'
$arr[3,4] = $null
$arr[3,5] = 'NO'
$arr[4,0] = '3736'
$arr[4,1] = 7
$arr[4,2] = '2'
$arr[4,3] = '0'
$arr[4,4] = '13.11%'
$arr[4,5] = 'NO'
$arr[5,0] = '3740'
$arr[5,1] = $null
$arr[5,2] = $null
$arr[5,3] = $null
$arr[5,4] = $null
$arr[5,5] = 'NO'
$arr[6,0] = '3742'
$arr[6,1] = 6
$arr[6,2] = '4'
$arr[6,3] = '1'
$arr[6,4] = '10.43%'
$arr[6,5] = 'NO'
$arr[7,0] = '3752'
$arr[7,1] = 8
$arr[7,2] = '3'
$arr[7,3] = '1'
$arr[7,4] = '8.88%'
$arr[7,5] = 'NO'
$arr[8,0] = '3757'
$arr[8,1] = 8
$arr[8,2] = '3'
$arr[8,3] = '4'
$arr[8,4] = '13.55%'
$arr[8,5] = 'YES'
$arr[9,0] = '3762'
$arr[9,1] = 8
$arr[9,2] = $null
$arr[9,3] = $null
$arr[9,4] = $null
$arr[9,5] = 'NO'
$arr[10,0] = '3766'
$arr[10,1] = 8
$arr[10,2] = '0'
$arr[10,3] = '0'
$arr[10,4] = $null
$arr[10,5] = 'NO'
$arr[11,0] = '3775'
$arr[11,1] = 8
$arr[11,2] = '0'
$arr[11,3] = '1'
$arr[11,4] = '4.40%'
$arr[11,5] = 'NO'
$arr[12,0] = '3788'
$arr[12,1] = $null
$arr[12,2] = $null
$arr[12,3] = $null
$arr[12,4] = $null
$arr[12,5] = 'NO'
$arr[13,0] = '3793'
$arr[13,1] = 8
$arr[13,2] = '1'
$arr[13,3] = '2'
$arr[13,4] = '8.00%'
$arr[13,5] = 'NO'
$arr[14,0] = '3852'
$arr[14,1] = $null
$arr[14,2] = $null
$arr[14,3] = $null
$arr[14,4] = $null
$arr[14,5] = 'NO'
$arr[15,0] = '4179'
$arr[15,1] = $null
$arr[15,2] = $null
$arr[15,3] = $null
$arr[15,4] = $null
$arr[15,5] = 'NO'
$arr[16,0] = '4304'
$arr[16,1] = 6
$arr[16,2] = $null
$arr[16,3] = $null
$arr[16,4] = $null
$arr[16,5] = 'NO'
$arr[17,0] = '4312'
$arr[17,1] = 7
$arr[17,2] = '2'
$arr[17,3] = '1'
$arr[17,4] = '5.49%'
$arr[17,5] = 'NO'
$arr[18,0] = '4321'
$arr[18,1] = 7
$arr[18,2] = '1'
$arr[18,3] = '2'
$arr[18,4] = '9.91%'
$arr[18,5] = 'NO'
$arr[19,0] = '4325'
$arr[19,1] = $null
$arr[19,2] = $null
$arr[19,3] = $null
$arr[19,4] = $null
$arr[19,5] = 'NO'
$extraSheet.Range($extraSheet.Cells.Item(2,1), $extraSheet.Cells.Item($nRows,$nCols)).Value = $arr

Write-Host "Workbook restructure complete."
